# Updated stats for Apr 8
# I42 was a projected/forecast value (formula extrapolating from prior
# days' growth rates). New actual data for 4/8 came in, so it becomes a
# literal "actual" value like the days before it (I41 and earlier),
# picking up the "actual data" cell formatting (style copied from I41)
# instead of the "projected data" formatting it had as a formula cell.
# Every dependent formula in row 42 (J/K/L/M/N) and the whole projection
# chain in I43:I49 (and their dependents) recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy I41's number formatting / fill (the "actual" style) onto I42 before
# overwriting its formula with a literal value.
$ws.Range("I41").Copy() | Out-Null
$ws.Range("I42").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Replace the forecast formula with the newly reported actual case count.
$ws.Range("I42").Value2 = 434927

# Move the active selection from I41 to I42, matching where the user was
# last working.
$ws.Range("I42").Select() | Out-Null
